$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
try {
  $ws1.Cells.Validation.Delete()
  Write-Host "all validation deleted"
} catch {
  Write-Host "ERROR: $_"
}
